$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert the new "if I refreshed the page..." bullet right before "Bugs:".
#    It goes directly after the "I think showing more relevant results..."
#    paragraph and inherits that paragraph's ListParagraph / numId=3 style.
# ---------------------------------------------------------------------------
$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "I think showing more relevant results*") {
        $anchorPara = $p
        break
    }
}
$anchorPara.Range.InsertParagraphAfter()

$newPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "I think showing more relevant results*") {
        $newPara = $d.Paragraphs.Item($i + 1)
        break
    }
}
$newPara.Range.Text = "if I refreshed the page, my results disappeared and I was unable to go back to my previous results."

# ---------------------------------------------------------------------------
# 2. Remove the existing "_GoBack" bookmark (currently wraps "My search ...
#    (??) " inside the "dog bones" bug paragraph).
# ---------------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 3. Merge the two runs that used to be split by that bookmark's end marker
#    ("helton song ... (??) " + "that would have to be ... products.") into a
#    single run, without disturbing any of the neighbouring runs. Wrapping
#    just that span in a temporary bookmark keeps the engine's run-merge
#    confined to this sub-range when we make (and immediately undo) a tiny
#    edit at the join point.
# ---------------------------------------------------------------------------
$fullText = $d.Content.Text
$mergeStart = $fullText.IndexOf("helton song with the word dog in the title")
$productsAt = $fullText.IndexOf("products.", $mergeStart)
$mergeEnd = $productsAt + "products.".Length

$mergeRange = $d.Range($mergeStart, $mergeEnd)
$d.Bookmarks.Add("_TmpMergeJoin", $mergeRange)

$joinPoint = $d.Range($mergeEnd, $mergeEnd)
$joinPoint.InsertBefore("X")

$fullText2 = $d.Content.Text
$stubAt = $fullText2.IndexOf("products.X")
$stubPos = $stubAt + "products.".Length
$stubRange = $d.Range($stubPos, $stubPos + 1)
$stubRange.Delete()

$d.Bookmarks.Item("_TmpMergeJoin").Delete()

# ---------------------------------------------------------------------------
# 4. Re-create "_GoBack" as an empty bookmark right at the start of the
#    "Bugs:" paragraph (before its only run).
# ---------------------------------------------------------------------------
$bugsPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Bugs:`r") {
        $bugsPara = $p
        break
    }
}
$bugsStart = $bugsPara.Range.Start
$goBackRange = $d.Range($bugsStart, $bugsStart)
$d.Bookmarks.Add("_GoBack", $goBackRange)
